$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 847, pushing existing rows 847:966 down to 849:968
$ws.Rows("847:848").Insert()

# --- Fill new row 847 (Primera) ---
$ws.Range("A847").Value = 6
$ws.Range("B847").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C847").Value = "Metropolitana"
$ws.Range("D847").Value = 44984
$ws.Range("E847").Value = 13
$ws.Range("F847").Value = 100112017
$ws.Range("G847").Value = "Apio"
$ws.Range("H847").Value = "Americana (o)"
$ws.Range("I847").Value = "Primera"
$ws.Range("J847").Value = 1040
$ws.Range("K847").Value = 7500
$ws.Range("L847").Value = 8000
$ws.Range("M847").Value = 7769
$ws.Range("N847").Value = "`$/docena de matas"
$ws.Range("O847").Value = "Región de Coquimbo"
$ws.Range("P847").Value = 1295
$ws.Range("Q847").Value = 6
$ws.Range("R847").Value = "Hortaliza"

# --- Fill new row 848 (Segunda) ---
$ws.Range("A848").Value = 6
$ws.Range("B848").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C848").Value = "Metropolitana"
$ws.Range("D848").Value = 44984
$ws.Range("E848").Value = 13
$ws.Range("F848").Value = 100112017
$ws.Range("G848").Value = "Apio"
$ws.Range("H848").Value = "Americana (o)"
$ws.Range("I848").Value = "Segunda"
$ws.Range("J848").Value = 290
$ws.Range("K848").Value = 6000
$ws.Range("L848").Value = 6000
$ws.Range("M848").Value = 6000
$ws.Range("N848").Value = "`$/docena de matas"
$ws.Range("O848").Value = "Región de Coquimbo"
$ws.Range("P848").Value = 1000
$ws.Range("Q848").Value = 6
$ws.Range("R848").Value = "Hortaliza"
